# "updated logbook + backlog"
# The backlog item "Permettre a n'importe quelle utilisateur de pouvoir
# noter une histoire" (row 2, column D - Statut) moves from "A faire" to
# "Fait". This introduces a brand-new shared string ("Fait") and leaves
# the cursor/selection sitting on D2, the cell that was just edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Fait"

# Leave the active selection on the cell that was just changed.
$ws.Range("D2").Select()
